$d = $word.ActiveDocument

# P11R0
$found = $d.Content.Find.Execute("Tubos e Tubulações - Definições`vTubos: Materiais, Processos de Fabricação e Normalização Dimensional.`vMeios de Ligação de Tubos, Conexões de Tubulações e Juntas de Expansão.`vVálvulas`vPurgadores de Vapor, Separadores e Filtros. Recomendações de Material para Serviços. `vAquecimento, Isolamento Térmico, Pintura e Proteção.`vDisposição das Construções em uma Instalação Industrial. Arranjo e Detalhamento de Tubulações.`vSistemas Especiais de Tubulação. Suportes de Tubulação. Montagem e Teste de Tubulações.`vVisita Técnica `vDesenhos de Tubulações`vDesenhos de Tubulações - Exercícios`vCálculo do diâmetro das tubulações`vA Tubulação Considerada como Elemento Estrutural Cálculo da Espessura de Parede de Tubos e do Vão entre Suportes.`vDilatação Térmica e Flexibilidade de Tubulações. Cálculo de Flexibilidade.`vCálculo de Flexibilidade.", $true, $false, $false, $false, $false, $true, 1, $false, "Provas em sala, entrega de exercícios ou casos práticos elaborados fora de sala de aula.", 2)
if (-not $found) { Write-Output "FAILED: P11R0" }

# P9R0
$found = $d.Content.Find.Execute("Tubos e Tubulações - Definições`vTubos: Materiais, Processos de Fabricação e Normalização Dimensional.`vMeios de Ligação de Tubos, Conexões de Tubulações e Juntas de Expansão.`vVálvulas`vPurgadores de Vapor, Separadores e Filtros. Recomendações de Material para Serviços. `vAquecimento, Isolamento Térmico, Pintura e Proteção.`vDisposição das Construções em uma Instalação Industrial. Arranjo e Detalhamento de Tubulações.`vSistemas Especiais de Tubulação. Suportes de Tubulação. Montagem e Teste de Tubulações.`vVisita Técnica `vDesenhos de Tubulações`vCálculo do diâmetro das tubulações`vA Tubulação Considerada como Elemento Estrutural Cálculo da Espessura de Parede de Tubos e do Vão entre Suportes.`vDilatação Térmica e Flexibilidade de Tubulações. Cálculo de `vCálculo de Flexibilidade.", $true, $false, $false, $false, $false, $true, 1, $false, "Aulas expositivas, desenvolvimento de exercícios em sala e fora de sala de aula. discussão de castos práticos, visitas técnicas", 2)
if (-not $found) { Write-Output "FAILED: P9R0" }

# P13R5
$found = $d.Content.Find.Execute("Frequência mínima de 70% e nota igual ou superior a 3,00 e inferior a 5,00 possibilita prova de recuperação.", $true, $false, $false, $false, $false, $true, 1, $false, "6634418 - Antonio Clelio Ribeiro", 2)
if (-not $found) { Write-Output "FAILED: P13R5" }

# P7R1
$found = $d.Content.Find.Execute("1285870 - Marcos Villela Barcza", $true, $false, $false, $false, $false, $true, 1, $false, "Tubos e Tubulações - Definições`vTubos: Materiais, Processos de Fabricação e Normalização Dimensional.`vMeios de Ligação de Tubos, Conexões de Tubulações e Juntas de Expansão.`vVálvulas`vPurgadores de Vapor, Separadores e Filtros. Recomendações de Material para Serviços. `vAquecimento, Isolamento Térmico, Pintura e Proteção.`vDisposição das Construções em uma Instalação Industrial. Arranjo e Detalhamento de Tubulações.`vSistemas Especiais de Tubulação. Suportes de Tubulação. Montagem e Teste de Tubulações.`vVisita Técnica `vDesenhos de Tubulações`vDesenhos de Tubulações - Exercícios`vCálculo do diâmetro das tubulações`vA Tubulação Considerada como Elemento Estrutural Cálculo da Espessura de Parede de Tubos e do Vão entre Suportes.`vDilatação Térmica e Flexibilidade de Tubulações. Cálculo de Flexibilidade.`vCálculo de Flexibilidade.", 2)
if (-not $found) { Write-Output "FAILED: P7R1" }

# P5R0
$found = $d.Content.Find.Execute("Ensinar a identificação e especificação dos elementos que compõem as tubulações que integram os processos inerentes às industrias de processamento.`vAuxiliar o desenvolvimento da habilidade de planejamento e projeto de processos industriais.", $true, $false, $false, $false, $false, $true, 1, $false, "Tubos e Tubulações - Definições`vTubos: Materiais, Processos de Fabricação e Normalização Dimensional.`vMeios de Ligação de Tubos, Conexões de Tubulações e Juntas de Expansão.`vVálvulas`vPurgadores de Vapor, Separadores e Filtros. Recomendações de Material para Serviços. `vAquecimento, Isolamento Térmico, Pintura e Proteção.`vDisposição das Construções em uma Instalação Industrial. Arranjo e Detalhamento de Tubulações.`vSistemas Especiais de Tubulação. Suportes de Tubulação. Montagem e Teste de Tubulações.`vVisita Técnica `vDesenhos de Tubulações`vCálculo do diâmetro das tubulações`vA Tubulação Considerada como Elemento Estrutural Cálculo da Espessura de Parede de Tubos e do Vão entre Suportes.`vDilatação Térmica e Flexibilidade de Tubulações. Cálculo de `vCálculo de Flexibilidade.", 2)
if (-not $found) { Write-Output "FAILED: P5R0" }

# P13R1
$found = $d.Content.Find.Execute("Aulas expositivas, desenvolvimento de exercícios em sala e fora de sala de aula. discussão de castos práticos, visitas técnicas`v", $true, $false, $false, $false, $false, $true, 1, $false, "Frequência mínima de 70% e nota igual ou superior a 3,00 e inferior a 5,00 possibilita prova de recuperação.`v", 2)
if (-not $found) { Write-Output "FAILED: P13R1" }

# P15R0
$found = $d.Content.Find.Execute("1)TUBULAÇÕES INDUSTRIAIS - Volume I e II`vSilva Telles, Pedro c. - Ed. Livros Técnicos e Científicos Editora S/A`v2)TABELAS E GRÁFICOS PARA PROJETOS DE TUBULAÇÕES INDUSTRIAIS`vSilva Telles, P.C./Paula Barros, Darcy G. - Ed. Interciência Ltda`v3)TUBULAÇÕES`vSilva, Remi Benedito - Editora Grêmio Politécnico da USP`v4)MATERIAIS PARA EQUIPAMENTOS DE PROCESSO`vSilva Telles, Pedro C. - Ed. Interciência Ltda`v5)CATÁLOGOS DIVERSOS", $true, $false, $false, $false, $false, $true, 1, $false, "1285870 - Marcos Villela Barcza", 2)
if (-not $found) { Write-Output "FAILED: P15R0" }

# P7R0
$found = $d.Content.Find.Execute("6634418 - Antonio Clelio Ribeiro`v", $true, $false, $false, $false, $false, $true, 1, $false, "Ensinar a identificação e especificação dos elementos que compõem as tubulações que integram os processos inerentes às industrias de processamento.`vAuxiliar o desenvolvimento da habilidade de planejamento e projeto de processos industriais.`v", 2)
if (-not $found) { Write-Output "FAILED: P7R0" }

# P13R3
$found = $d.Content.Find.Execute("Provas em sala, entrega de exercícios ou casos práticos elaborados fora de sala de aula.`v", $true, $false, $false, $false, $false, $true, 1, $false, "1)TUBULAÇÕES INDUSTRIAIS - Volume I e II`vSilva Telles, Pedro c. - Ed. Livros Técnicos e Científicos Editora S/A`v2)TABELAS E GRÁFICOS PARA PROJETOS DE TUBULAÇÕES INDUSTRIAIS`vSilva Telles, P.C./Paula Barros, Darcy G. - Ed. Interciência Ltda`v3)TUBULAÇÕES`vSilva, Remi Benedito - Editora Grêmio Politécnico da USP`v4)MATERIAIS PARA EQUIPAMENTOS DE PROCESSO`vSilva Telles, Pedro C. - Ed. Interciência Ltda`v5)CATÁLOGOS DIVERSOS`v", 2)
if (-not $found) { Write-Output "FAILED: P13R3" }
